# Auto-generated edit script: updates Leve profit/price data cells across
# the ALC, ARM, BSM, CRP, CUL, GSM, LTW, and WVR sheets to match refreshed
# market data from the scheduled scraping run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15 (ALC)
$ws.Range("H15").Value = 1273.091
$ws.Range("I15").Value = 1273.091
$ws.Range("K15").Value = 3819.273
$ws.Range("M15").Value = -3650.273

# Row 62 (ALC)
$ws.Range("H62").Value = 9339.315000000001
$ws.Range("I62").Value = 8430.875
$ws.Range("K62").Value = 8430.875
$ws.Range("M62").Value = -7806.875

# Row 65 (ALC)
$ws.Range("H65").Value = 9339.315000000001
$ws.Range("I65").Value = 8430.875
$ws.Range("K65").Value = 42154.375
$ws.Range("M65").Value = -39034.375

# Row 125 (ALC)
$ws.Range("H125").Value = 1351.375
$ws.Range("J125").Value = 2000
$ws.Range("L125").Value = 18000
$ws.Range("N125").Value = -22920

# Row 137 (ALC)
$ws.Range("H137").Value = 8886.166999999999
$ws.Range("I137").Value = 1665
$ws.Range("J137").Value = 10949.357
$ws.Range("K137").Value = 4995
$ws.Range("L137").Value = 32848.071
$ws.Range("M137").Value = -2445
$ws.Range("N137").Value = -37948.071

# Row 138 (ALC)
$ws.Range("H138").Value = 3298.6394
$ws.Range("J138").Value = 4268.0513
$ws.Range("L138").Value = 12804.1539
$ws.Range("N138").Value = -23084.1539

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (ARM)
$ws.Range("H2").Value = 1877.65
$ws.Range("I2").Value = 1927.625
$ws.Range("K2").Value = 1927.625
$ws.Range("M2").Value = -1814.625

# Row 32 (ARM)
$ws.Range("H32").Value = 2960.6738
$ws.Range("I32").Value = 2534.0264
$ws.Range("K32").Value = 2534.0264
$ws.Range("M32").Value = -2247.0264

# Row 61 (ARM)
$ws.Range("H61").Value = 5000
$ws.Range("I61").Value = 5000
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 5000
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -4788
$ws.Range("N61").ClearContents()

# Row 74 (ARM)
$ws.Range("H74").Value = 12776.1
$ws.Range("I74").Value = 3095.125
$ws.Range("K74").Value = 3095.125
$ws.Range("M74").Value = -2221.125

# Row 77 (ARM)
$ws.Range("H77").Value = 12776.1
$ws.Range("I77").Value = 3095.125
$ws.Range("K77").Value = 15475.625
$ws.Range("M77").Value = -11107.625

# Row 102 (ARM)
$ws.Range("H102").Value = 3132
$ws.Range("I102").Value = 2997.5173
$ws.Range("K102").Value = 2997.5173
$ws.Range("M102").Value = -1375.5173

# Row 110 (ARM)
$ws.Range("H110").Value = 9769.929
$ws.Range("I110").Value = 15342.308
$ws.Range("K110").Value = 15342.308
$ws.Range("M110").Value = -13297.308

# Row 116 (ARM)
$ws.Range("H116").Value = 1877.65
$ws.Range("I116").Value = 1927.625
$ws.Range("K116").Value = 1927.625
$ws.Range("M116").Value = 366.375

# Row 122 (ARM)
$ws.Range("H122").Value = 1809.9445
$ws.Range("I122").Value = 1936.5625
$ws.Range("J122").Value = 797
$ws.Range("K122").Value = 5809.6875
$ws.Range("L122").Value = 2391
$ws.Range("M122").Value = -3359.6875
$ws.Range("N122").Value = -7291

# Row 132 (ARM)
$ws.Range("H132").Value = 3003.4736
$ws.Range("I132").Value = 2365.5386
$ws.Range("J132").Value = 4385.6665
$ws.Range("K132").Value = 7096.6158
$ws.Range("L132").Value = 13156.9995
$ws.Range("M132").Value = -4566.6158
$ws.Range("N132").Value = -18216.9995

# Row 136 (ARM)
$ws.Range("H136").Value = 5000
$ws.Range("I136").Value = 5000
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 15000
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -12450
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (BSM)
$ws.Range("H3").Value = 1877.65
$ws.Range("I3").Value = 1927.625
$ws.Range("K3").Value = 1927.625
$ws.Range("M3").Value = -1813.625

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 9311186
$ws.Range("I31").Value = 18588596
$ws.Range("K31").Value = 18588596
$ws.Range("M31").Value = -18588301

# Row 34 (CRP)
$ws.Range("H34").Value = 9311186
$ws.Range("I34").Value = 18588596
$ws.Range("K34").Value = 18588596
$ws.Range("M34").Value = -18588394

# Row 86 (CRP)
$ws.Range("H86").Value = 4826
$ws.Range("I86").Value = 4741.4
$ws.Range("K86").Value = 4741.4
$ws.Range("M86").Value = -3618.4

# Row 89 (CRP)
$ws.Range("H89").Value = 4826
$ws.Range("I89").Value = 4741.4
$ws.Range("K89").Value = 23707
$ws.Range("M89").Value = -18091

# Row 99 (CRP)
$ws.Range("H99").Value = 3638.125
$ws.Range("I99").Value = 3534.125
$ws.Range("J99").Value = 3846.125
$ws.Range("K99").Value = 3534.125
$ws.Range("L99").Value = 3846.125
$ws.Range("M99").Value = -2036.125
$ws.Range("N99").Value = -6842.125

# Row 126 (CRP)
$ws.Range("H126").Value = 3638.125
$ws.Range("I126").Value = 3534.125
$ws.Range("J126").Value = 3846.125
$ws.Range("K126").Value = 10602.375
$ws.Range("L126").Value = 11538.375
$ws.Range("M126").Value = -8132.375
$ws.Range("N126").Value = -16478.375

# Row 132 (CRP)
$ws.Range("H132").Value = 2619
$ws.Range("I132").Value = 2306.2
$ws.Range("J132").Value = 3661.6667
$ws.Range("K132").Value = 6918.599999999999
$ws.Range("L132").Value = 10985.0001
$ws.Range("M132").Value = -4388.599999999999
$ws.Range("N132").Value = -16045.0001

$ws = $wb.Worksheets.Item("CUL")
# Row 68 (CUL)
$ws.Range("H68").Value = 5559188
$ws.Range("I68").Value = 4596.5
$ws.Range("J68").Value = 7146214
$ws.Range("K68").Value = 13789.5
$ws.Range("L68").Value = 21438642
$ws.Range("M68").Value = -12978.5
$ws.Range("N68").Value = -21440264

# Row 71 (CUL)
$ws.Range("H71").Value = 5559188
$ws.Range("I71").Value = 4596.5
$ws.Range("J71").Value = 7146214
$ws.Range("K71").Value = 41368.5
$ws.Range("L71").Value = 64315926
$ws.Range("M71").Value = -37312.5
$ws.Range("N71").Value = -64324038

# Row 92 (CUL)
$ws.Range("H92").Value = 565.2857
$ws.Range("I92").Value = 244
$ws.Range("J92").Value = 806.25
$ws.Range("K92").Value = 732
$ws.Range("L92").Value = 2418.75
$ws.Range("M92").Value = 516
$ws.Range("N92").Value = -4914.75

$ws = $wb.Worksheets.Item("GSM")
# Row 44 (GSM)
$ws.Range("H44").Value = 12676
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

# Row 52 (GSM)
$ws.Range("H52").Value = 22000
$ws.Range("J52").Value = 22000
$ws.Range("L52").Value = 22000
$ws.Range("N52").Value = -22518

# Row 70 (GSM)
$ws.Range("H70").Value = 9370
$ws.Range("I70").Value = 11281.167
$ws.Range("K70").Value = 11281.167
$ws.Range("M70").Value = -11011.167

# Row 73 (GSM)
$ws.Range("H73").Value = 9370
$ws.Range("I73").Value = 11281.167
$ws.Range("K73").Value = 11281.167
$ws.Range("M73").Value = -10345.167

# Row 104 (GSM)
$ws.Range("H104").Value = 140000
$ws.Range("J104").Value = 140000
$ws.Range("L104").Value = 140000
$ws.Range("N104").Value = -146988

# Row 113 (GSM)
$ws.Range("H113").Value = 3323.4194
$ws.Range("I113").Value = 2824.5557
$ws.Range("K113").Value = 2824.5557
$ws.Range("M113").Value = -654.5556999999999

# Row 122 (GSM)
$ws.Range("H122").Value = 532.875
$ws.Range("I122").Value = 514.6667
$ws.Range("K122").Value = 1544.0001
$ws.Range("M122").Value = 905.9999

# Row 132 (GSM)
$ws.Range("H132").Value = 246889.14
$ws.Range("I132").Value = 280485.06
$ws.Range("J132").Value = 4998.6
$ws.Range("K132").Value = 841455.1799999999
$ws.Range("L132").Value = 14995.8
$ws.Range("M132").Value = -838925.1799999999
$ws.Range("N132").Value = -20055.8

$ws = $wb.Worksheets.Item("LTW")
# Row 40 (LTW)
$ws.Range("H40").Value = 5321.375
$ws.Range("I40").Value = 5335.1304
$ws.Range("K40").Value = 5335.1304
$ws.Range("M40").Value = -5199.1304

# Row 46 (LTW)
$ws.Range("H46").Value = 3956.7144
$ws.Range("I46").Value = 1939.6
$ws.Range("K46").Value = 1939.6
$ws.Range("M46").Value = -1751.6

# Row 68 (LTW)
$ws.Range("H68").Value = 4999.4614
$ws.Range("J68").Value = 8664.666999999999
$ws.Range("L68").Value = 8664.666999999999
$ws.Range("N68").Value = -10162.667

# Row 71 (LTW)
$ws.Range("H71").Value = 4999.4614
$ws.Range("J71").Value = 8664.666999999999
$ws.Range("L71").Value = 43323.335
$ws.Range("N71").Value = -50811.335

# Row 132 (LTW)
$ws.Range("H132").Value = 5145.619
$ws.Range("I132").Value = 4289.857
$ws.Range("J132").Value = 6857.143
$ws.Range("K132").Value = 12869.571
$ws.Range("L132").Value = 20571.429
$ws.Range("M132").Value = -10339.571
$ws.Range("N132").Value = -25631.429

$ws = $wb.Worksheets.Item("WVR")
# Row 136 (WVR)
$ws.Range("H136").Value = 2598.75
$ws.Range("I136").Value = 2190.7693
$ws.Range("K136").Value = 6572.3079
$ws.Range("M136").Value = -4022.3079
